$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at O and P (shifting old O:U to Q:W)
$ws.Range("O1:P1").EntireColumn.Insert()

# Update header labels affected by the rename/insert
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Copy style of neighboring header cell (N1) onto the two new header cells
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the numeric values in columns M, N, O, P for each data row
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 3
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4

$ws.Range("M4").Value = 3
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 3

$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 2

$ws.Range("M6").Value = 2
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
